$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = 959
$ws.Range("D8").Value = 161
$ws.Range("E8").Value = 798
$ws.Range("F8").Value = 6.60377358490566
$ws.Range("G8").Value = 83.21167883211679
$ws.Range("H8").Value = 16.78832116788321
